# edit.ps1 - applies the "Revise and update User Guide pages" changes
# to webapp/static/user_guide/importing_xml.docx
#
# Strategy: the document's paragraph count/order does not change (no
# paragraphs are inserted or removed by the target diff), so we can
# address each paragraph by its stable 1-based index and perform
# targeted Find/Replace (scoped to that paragraph's Range, replacing
# exactly one occurrence) plus a couple of direct Range edits for
# insertions and character formatting.

$d = $word.ActiveDocument

# wdReplaceOne = 1 ; wdFindContinue (not used, we always pass Wrap=1 for scope search)
$wdReplaceOne = 1

# ---------------------------------------------------------------
# 1) Paragraph 3: "...outside of ezEML." -> "...outside of ezEML. ezEML
#    can do this, but only up to a point."
# ---------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Find.Execute(
    "outside of ezEML.", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "outside of ezEML. ezEML can do this, but only up to a point.",
    $wdReplaceOne) | Out-Null

# ---------------------------------------------------------------
# 2) Paragraph 5: several wording tweaks + new leading "I.e., "
# ---------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)

$p5.Range.Find.Execute(
    "ezEML program can", $false, $false, $false, $false, $false,
    $true, 1, $false, "ezEML software can", $wdReplaceOne) | Out-Null

$p5.Range.Find.Execute(
    "elements it encounters are", $false, $false, $false, $false, $false,
    $true, 1, $false, "elements it will encounter are", $wdReplaceOne) | Out-Null

$p5.Range.Find.Execute(
    "ezEML has to drop", $false, $false, $false, $false, $false,
    $true, 1, $false, "ezEML needs to drop", $wdReplaceOne) | Out-Null

$p5.Range.Find.Execute(
    "complex than in the pure", $false, $false, $false, $false, $false,
    $true, 1, $false, "complex than it would be in the pure", $wdReplaceOne) | Out-Null

# Prepend "I.e., " at the very start of the paragraph.
$p5Start = $p5.Range.Start
$d.Range($p5Start, $p5Start).InsertBefore("I.e., ")

# ---------------------------------------------------------------
# 3) Mark both screenshots as NoProof (adds <w:rPr><w:noProof/></w:rPr>
#    to the run that hosts each <w:drawing>).
# ---------------------------------------------------------------
$d.InlineShapes.Item(1).Range.Font.NoProofing = 1
$d.InlineShapes.Item(2).Range.Font.NoProofing = 1

# ---------------------------------------------------------------
# 4) Paragraph 9 (the "Select an XML file and click Import." paragraph):
#    - "the subset of" -> "the features of"
#    - "...scope, ezEML will" -> "...scope, however, ezEML will"
#    - "following example:" -> "following screenshot showing the result
#       of a sample import:"
# ---------------------------------------------------------------
$p9 = $d.Paragraphs.Item(9)
$rightQuote = [char]0x2019

$p9.Range.Find.Execute(
    "the subset of", $false, $false, $false, $false, $false,
    $true, 1, $false, "the features of", $wdReplaceOne) | Out-Null

$p9.Range.Find.Execute(
    ("ezEML" + $rightQuote + "s scope, ezEML will"), $false, $false, $false, $false, $false,
    $true, 1, $false,
    ("ezEML" + $rightQuote + "s scope, however, ezEML will"), $wdReplaceOne) | Out-Null

$p9.Range.Find.Execute(
    "following example:", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "following screenshot showing the result of a sample import:",
    $wdReplaceOne) | Out-Null

# ---------------------------------------------------------------
# 5) Paragraph 17: italicize the word "text" that immediately follows
#    "include" (but not the later "such text elements" occurrence).
# ---------------------------------------------------------------
$p17 = $d.Paragraphs.Item(17)
$p17Text = $p17.Range.Text
$localIdx = $p17Text.IndexOf("include text elements")
$wordStart = $p17.Range.Start + $localIdx + ("include ").Length
$wordEnd = $wordStart + ("text").Length
$d.Range($wordStart, $wordEnd).Font.Italic = 1

Write-Host "Done."
